# Evalresults.xlsx -- add participant 12's second (RecetteTek) SUS answers,
# their comments/preference, extend the shared SUS-score formulas down to
# row 15, and refresh the two Average cells in row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: RecetteTek SUS answers (C:L) for participant 12 -----------
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 4
$ws.Range("L15").Value = 3

# --- Row 15: Recipes For Life SUS answers (P:Y) for participant 12 -----
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 4
$ws.Range("S15").Value = 4
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 4
$ws.Range("V15").Value = 4
$ws.Range("W15").Value = 3
$ws.Range("X15").Value = 3
$ws.Range("Y15").Value = 2

# --- Extend the SUS-total formulas (columns M and Z) down to row 15 ----
# Re-applying the formula across M6:M15 / Z6:Z15 turns them back into a
# single shared-formula group covering the new row, same as Excel does
# when you drag-fill the existing formula down one row.
$ws.Range("M6:M15").Formula = "=( (C6 - 1) + (E6 - 1) + (G6 - 1) + (I6-1) + (K6-1)  + (5 -D6) + (5 - F6) + (5 - H6) + (5 - J6) + (5-L6))*2.5"
$ws.Range("Z6:Z15").Formula = "=( (P6 - 1) + (R6 - 1) + (T6 - 1) + (V6-1) + (X6-1)  + (5 -Q6) + (5 - S6) + (5 - U6) + (5 - W6) + (5-Y6))*2.5"

# --- Comments / application preference for participant 12 --------------
$ws.Range("AA15").Value = "It was mostly straightfoward to use"
$ws.Range("AA15").WrapText = $true
$ws.Range("AC15").Value = "Recipes For Life cause you could do more. The first one was slightly easier to use"
$ws.Range("AC15").WrapText = $true

# Row 15 grows to fit the new wrapped comment text.
$ws.Rows(15).RowHeight = 72.6

# --- Row 17 averages now span the 10 participants (rows 6:15) ----------
$ws.Range("M17").Formula = "=SUM(M6:M15)  / 10"
$ws.Range("Z17").Formula = "=SUM(Z6:Z15)/10"

# --- Scroll the sheet view down and select the new total cell ----------
$ws.Activate()
$ws.Range("M15").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
